$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# All edited cells live in D2:E50. Toggle that block to a text format
# before writing so values keep their exact original formatting
# (trailing zeros, "%" suffix, etc.) instead of being parsed as numbers,
# then restore the default style afterwards.
$fmtRange = $ws.Range("D2:E50")
$fmtRange.NumberFormat = "@"

$ws.Range("D2").Value = "278.68"
$ws.Range("E2").Value = "6.65%"
$ws.Range("D3").Value = "27.13"
$ws.Range("D4").Value = "4.835"
$ws.Range("E4").Value = "2.92%"
$ws.Range("D5").Value = "0.06260"
$ws.Range("E5").Value = "0.70%"
$ws.Range("D6").Value = "6.858"
$ws.Range("E6").Value = "1.61%"
$ws.Range("D7").Value = "0.8787"
$ws.Range("E7").Value = "3.14%"
$ws.Range("D8").Value = "0.9403"
$ws.Range("E8").Value = "2.43%"
$ws.Range("D9").Value = "0.1451"
$ws.Range("E9").Value = "3.57%"
$ws.Range("D10").Value = "0.05147"
$ws.Range("E10").Value = "6.29%"
$ws.Range("D11").Value = "0.07287"
$ws.Range("E11").Value = "2.80%"
$ws.Range("D12").Value = "0.03162"
$ws.Range("E12").Value = "1.71%"
$ws.Range("D13").Value = "0.09045"
$ws.Range("E13").Value = "-0.11%"
$ws.Range("E14").Value = "0.86%"
$ws.Range("D15").Value = "0.0006277"
$ws.Range("E15").Value = "1.87%"
$ws.Range("D16").Value = "0.006071"
$ws.Range("E16").Value = "-0.59%"
$ws.Range("D17").Value = "3.451"
$ws.Range("E17").Value = "0.29%"
$ws.Range("D18").Value = "3.264"
$ws.Range("E18").Value = "3.02%"
$ws.Range("D19").Value = "2.286"
$ws.Range("E19").Value = "5.60%"
$ws.Range("E20").Value = "-0.62%"
$ws.Range("E21").Value = "0.05%"
$ws.Range("D22").Value = "3.850"
$ws.Range("E22").Value = "-5.93%"
$ws.Range("D23").Value = "0.04315"
$ws.Range("E23").Value = "1.45%"
$ws.Range("D24").Value = "0.001177"
$ws.Range("D25").Value = "0.004274"
$ws.Range("E25").Value = "4.77%"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").Value = "-0.05%"
$ws.Range("D27").Value = "0.0001689"
$ws.Range("E27").Value = "3.01%"
$ws.Range("D40").Value = "0.04027"
$ws.Range("E40").Value = "1.55%"
$ws.Range("D41").Value = "0.006403"
$ws.Range("E41").Value = "55.74%"
$ws.Range("D42").Value = "0.1153"
$ws.Range("E42").Value = "3.68%"
$ws.Range("D43").Value = "0.002160"
$ws.Range("E43").Value = "-2.31%"
$ws.Range("D44").Value = "0.01398"
$ws.Range("E44").Value = "2.40%"
$ws.Range("D45").Value = "0.00005177"
$ws.Range("E45").Value = "0.29%"
$ws.Range("D46").Value = "0.00000000750"
$ws.Range("E46").Value = "-0.10%"
$ws.Range("D47").Value = "2.366"
$ws.Range("E47").Value = "755.33%"
$ws.Range("D49").Value = "0.00002099"
$ws.Range("E49").Value = "-0.10%"
$ws.Range("D50").Value = "0.0001999"
$ws.Range("E50").Value = "-0.10%"

$fmtRange.Style = "Normal"

